$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the sequence batch name in SB26-AB339: "SB26" -> "SB26_AB339"
$ws.Range("B2").Value = "SB26_AB339"

# Move the active selection to C5 (matches resulting sheetView selection)
$ws.Range("C5").Select()
